# Add gcamp_type and expt_group
# Add these two columns to Info and record in data.mat

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend formatting of the last existing column (R) into the new column (S)
# so the new cells pick up the same style index as their neighbors.
$ws.Range("R1:R2").AutoFill($ws.Range("R1:S2"), 0)

# Populate the new column: header "expt_group" and value "NDNF"
$ws.Range("S1").Value = "expt_group"
$ws.Range("S2").Value = "NDNF"

# Update the selected cell as recorded in the saved view
$ws.Range("D5").Select()
